# Auto-generated edit script applying the Diabolos_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 4
$ws.Range("H4").Value2 = 166668340
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 166668340
$ws.Range("K4").Value2 = 0
$ws.Range("L4").Value2 = 166668340
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value2 = -166668568
# ALC row 40
$ws.Range("H40").Value2 = 2535.2942
$ws.Range("J40").Value2 = 2700
$ws.Range("L40").Value2 = 2700
$ws.Range("N40").Value2 = -3050
# ALC row 76
$ws.Range("H76").Value2 = 9650981
$ws.Range("I76").Value2 = 130499.25
$ws.Range("J76").Value2 = 17267366
$ws.Range("K76").Value2 = 130499.25
$ws.Range("L76").Value2 = 17267366
$ws.Range("M76").Value2 = -130184.25
$ws.Range("N76").Value2 = -17267996
# ALC row 79
$ws.Range("H79").Value2 = 9650981
$ws.Range("I79").Value2 = 130499.25
$ws.Range("J79").Value2 = 17267366
$ws.Range("K79").Value2 = 130499.25
$ws.Range("L79").Value2 = 17267366
$ws.Range("M79").Value2 = -129407.25
$ws.Range("N79").Value2 = -17269550
# ALC row 107
$ws.Range("H107").Value2 = 587.13794
$ws.Range("I107").Value2 = 608.375
$ws.Range("K107").Value2 = 608.375
$ws.Range("M107").Value2 = 1311.625
# ALC row 132
$ws.Range("H132").Value2 = 4405.952
$ws.Range("I132").Value2 = 4215.943
$ws.Range("K132").Value2 = 12647.829
$ws.Range("M132").Value2 = -10117.829

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value2 = 5143.684
$ws.Range("I32").Value2 = 6520
$ws.Range("K32").Value2 = 6520
$ws.Range("M32").Value2 = -6233
# ARM row 132
$ws.Range("H132").Value2 = 1207.375
$ws.Range("I132").Value2 = 1002.3111
$ws.Range("J132").Value2 = 2046.2727
$ws.Range("K132").Value2 = 3006.9333
$ws.Range("L132").Value2 = 6138.8181
$ws.Range("M132").Value2 = -476.9333000000001
$ws.Range("N132").Value2 = -11198.8181

$ws = $wb.Worksheets.Item("BSM")
# BSM row 22
$ws.Range("H22").Value2 = 0
$ws.Range("I22").Value2 = 0
$ws.Range("K22").Value2 = 0
$ws.Range("M22").ClearContents()
# BSM row 80
$ws.Range("H80").Value2 = 25213.375
$ws.Range("J80").Value2 = 66817.336
$ws.Range("L80").Value2 = 66817.336
$ws.Range("N80").Value2 = -68813.336
# BSM row 83
$ws.Range("H83").Value2 = 25213.375
$ws.Range("J83").Value2 = 66817.336
$ws.Range("L83").Value2 = 334086.68
$ws.Range("N83").Value2 = -344070.68
# BSM row 86
$ws.Range("H86").Value2 = 27780118
$ws.Range("I86").Value2 = 31252340
$ws.Range("J86").Value2 = 2353
$ws.Range("K86").Value2 = 31252340
$ws.Range("L86").Value2 = 2353
$ws.Range("M86").Value2 = -31251217
$ws.Range("N86").Value2 = -4599
# BSM row 89
$ws.Range("H89").Value2 = 27780118
$ws.Range("I89").Value2 = 31252340
$ws.Range("J89").Value2 = 2353
$ws.Range("K89").Value2 = 156261700
$ws.Range("L89").Value2 = 11765
$ws.Range("M89").Value2 = -156256084
$ws.Range("N89").Value2 = -22997
# BSM row 94
$ws.Range("H94").Value2 = 4042
$ws.Range("I94").Value2 = 4800.933
$ws.Range("K94").Value2 = 4800.933
$ws.Range("M94").Value2 = -4349.933
# BSM row 107
$ws.Range("H107").Value2 = 1508.0667
$ws.Range("I107").Value2 = 1336.5238
$ws.Range("J107").Value2 = 1908.3334
$ws.Range("K107").Value2 = 1336.5238
$ws.Range("L107").Value2 = 1908.3334
$ws.Range("M107").Value2 = 583.4762000000001
$ws.Range("N107").Value2 = -5748.3334
# BSM row 134
$ws.Range("H134").Value2 = 971.0714
$ws.Range("I134").Value2 = 971.0714
$ws.Range("K134").Value2 = 2913.2142
$ws.Range("M134").Value2 = -378.2142000000003

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value2 = 5676.1523
$ws.Range("I31").Value2 = 16161.375
$ws.Range("K31").Value2 = 16161.375
$ws.Range("M31").Value2 = -15866.375
# CRP row 34
$ws.Range("H34").Value2 = 5676.1523
$ws.Range("I34").Value2 = 16161.375
$ws.Range("K34").Value2 = 16161.375
$ws.Range("M34").Value2 = -15959.375
# CRP row 58
$ws.Range("H58").Value2 = 2102.64
$ws.Range("I58").Value2 = 1855.6
$ws.Range("J58").Value2 = 2473.2
$ws.Range("K58").Value2 = 1855.6
$ws.Range("L58").Value2 = 2473.2
$ws.Range("M58").Value2 = -1652.6
$ws.Range("N58").Value2 = -2879.2
# CRP row 105
$ws.Range("H105").Value2 = 3442.111
$ws.Range("I105").Value2 = 2663.1667
$ws.Range("K105").Value2 = 2663.1667
$ws.Range("M105").Value2 = -916.1667000000002
# CRP row 122
$ws.Range("H122").Value2 = 2448.8635
$ws.Range("I122").Value2 = 2433.2942
$ws.Range("J122").Value2 = 2501.8
$ws.Range("K122").Value2 = 7299.882599999999
$ws.Range("L122").Value2 = 7505.400000000001
$ws.Range("M122").Value2 = -4849.882599999999
$ws.Range("N122").Value2 = -12405.4
# CRP row 132
$ws.Range("H132").Value2 = 1659.8695
$ws.Range("I132").Value2 = 1400.5555
$ws.Range("J132").Value2 = 2593.4
$ws.Range("K132").Value2 = 4201.666499999999
$ws.Range("L132").Value2 = 7780.200000000001
$ws.Range("M132").Value2 = -1671.666499999999
$ws.Range("N132").Value2 = -12840.2
# CRP row 134
$ws.Range("H134").Value2 = 2667.7856
$ws.Range("I134").Value2 = 2306.3635
$ws.Range("K134").Value2 = 6919.0905
$ws.Range("M134").Value2 = -4384.0905
# CRP row 136
$ws.Range("H136").Value2 = 2102.64
$ws.Range("I136").Value2 = 1855.6
$ws.Range("J136").Value2 = 2473.2
$ws.Range("K136").Value2 = 5566.799999999999
$ws.Range("L136").Value2 = 7419.599999999999
$ws.Range("M136").Value2 = -3016.799999999999
$ws.Range("N136").Value2 = -12519.6

$ws = $wb.Worksheets.Item("CUL")
# CUL row 98
$ws.Range("H98").Value2 = 818
$ws.Range("I98").Value2 = 726.5
$ws.Range("J98").Value2 = 1001
$ws.Range("K98").Value2 = 2179.5
$ws.Range("L98").Value2 = 3003
$ws.Range("M98").Value2 = -681.5
$ws.Range("N98").Value2 = -5999
# CUL row 113
$ws.Range("H113").Value2 = 1819.5
$ws.Range("I113").Value2 = 831.2
$ws.Range("J113").Value2 = 2525.4285
$ws.Range("K113").Value2 = 2493.6
$ws.Range("L113").Value2 = 7576.2855
$ws.Range("M113").Value2 = -323.6000000000004
$ws.Range("N113").Value2 = -11916.2855
# CUL row 123
$ws.Range("H123").Value2 = 2497.5
$ws.Range("I123").Value2 = 2497.5
$ws.Range("K123").Value2 = 7492.5
$ws.Range("M123").Value2 = -5042.5

$ws = $wb.Worksheets.Item("GSM")
# GSM row 25
$ws.Range("H25").Value2 = 2250
$ws.Range("J25").Value2 = 2250
$ws.Range("L25").Value2 = 2250
$ws.Range("N25").Value2 = -3308
# GSM row 102
$ws.Range("H102").Value2 = 1760.2727
$ws.Range("I102").Value2 = 1456.2778
$ws.Range("J102").Value2 = 3128.25
$ws.Range("K102").Value2 = 1456.2778
$ws.Range("L102").Value2 = 3128.25
$ws.Range("M102").Value2 = 165.7221999999999
$ws.Range("N102").Value2 = -6372.25
# GSM row 126
$ws.Range("H126").Value2 = 11375.917
$ws.Range("J126").Value2 = 3583.3333
$ws.Range("L126").Value2 = 10749.9999
$ws.Range("N126").Value2 = -15689.9999

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value2 = 3237.3333
$ws.Range("I7").Value2 = 2356.0833
$ws.Range("J7").Value2 = 4999.8335
$ws.Range("K7").Value2 = 2356.0833
$ws.Range("L7").Value2 = 4999.8335
$ws.Range("M7").Value2 = -2244.0833
$ws.Range("N7").Value2 = -5223.8335
# LTW row 16
$ws.Range("H16").Value2 = 578.4375
$ws.Range("I16").Value2 = 441.7857
$ws.Range("J16").Value2 = 1535
$ws.Range("K16").Value2 = 441.7857
$ws.Range("L16").Value2 = 1535
$ws.Range("M16").Value2 = -271.7857
$ws.Range("N16").Value2 = -1875
# LTW row 20
$ws.Range("H20").Value2 = 0
$ws.Range("J20").Value2 = 0
$ws.Range("L20").Value2 = 0
$ws.Range("N20").ClearContents()
# LTW row 22
$ws.Range("H22").Value2 = 2488219
$ws.Range("I22").Value2 = 1467.1428
$ws.Range("J22").Value2 = 5389429.5
$ws.Range("K22").Value2 = 1467.1428
$ws.Range("L22").Value2 = 5389429.5
$ws.Range("M22").Value2 = -1172.1428
$ws.Range("N22").Value2 = -5390019.5
# LTW row 27
$ws.Range("H27").Value2 = 2488219
$ws.Range("I27").Value2 = 1467.1428
$ws.Range("J27").Value2 = 5389429.5
$ws.Range("K27").Value2 = 1467.1428
$ws.Range("L27").Value2 = 5389429.5
$ws.Range("M27").Value2 = -1360.1428
$ws.Range("N27").Value2 = -5389643.5
# LTW row 55
$ws.Range("H55").Value2 = 411.9565
$ws.Range("I55").Value2 = 351.4375
$ws.Range("J55").Value2 = 550.2857
$ws.Range("K55").Value2 = 351.4375
$ws.Range("L55").Value2 = 550.2857
$ws.Range("M55").Value2 = -178.4375
$ws.Range("N55").Value2 = -896.2857
# LTW row 68
$ws.Range("H68").Value2 = 7025.4443
$ws.Range("I68").Value2 = 7141.4546
$ws.Range("J68").Value2 = 6843.143
$ws.Range("K68").Value2 = 7141.4546
$ws.Range("L68").Value2 = 6843.143
$ws.Range("M68").Value2 = -6392.4546
$ws.Range("N68").Value2 = -8341.143
# LTW row 71
$ws.Range("H71").Value2 = 7025.4443
$ws.Range("I71").Value2 = 7141.4546
$ws.Range("J71").Value2 = 6843.143
$ws.Range("K71").Value2 = 35707.273
$ws.Range("L71").Value2 = 34215.715
$ws.Range("M71").Value2 = -31963.273
$ws.Range("N71").Value2 = -41703.715
# LTW row 126
$ws.Range("H126").Value2 = 3237.3333
$ws.Range("I126").Value2 = 2356.0833
$ws.Range("J126").Value2 = 4999.8335
$ws.Range("K126").Value2 = 7068.249899999999
$ws.Range("L126").Value2 = 14999.5005
$ws.Range("M126").Value2 = -4598.249899999999
$ws.Range("N126").Value2 = -19939.5005
# LTW row 132
$ws.Range("H132").Value2 = 7005.85
$ws.Range("I132").Value2 = 3749.875
$ws.Range("J132").Value2 = 9176.5
$ws.Range("K132").Value2 = 11249.625
$ws.Range("L132").Value2 = 27529.5
$ws.Range("M132").Value2 = -8719.625
$ws.Range("N132").Value2 = -32589.5
# LTW row 136
$ws.Range("H136").Value2 = 2408.4243
$ws.Range("I136").Value2 = 2017.7778
$ws.Range("K136").Value2 = 6053.3334
$ws.Range("M136").Value2 = -3503.3334

$ws = $wb.Worksheets.Item("WVR")
# WVR row 30
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()
# WVR row 107
$ws.Range("H107").Value2 = 924.8570999999999
$ws.Range("I107").Value2 = 949.75
$ws.Range("K107").Value2 = 2849.25
$ws.Range("M107").Value2 = -929.25
# WVR row 122
$ws.Range("H122").Value2 = 2056.7856
$ws.Range("I122").Value2 = 1728.4286
$ws.Range("K122").Value2 = 5185.2858
$ws.Range("M122").Value2 = -2735.2858
# WVR row 126
$ws.Range("H126").Value2 = 1772.4546
$ws.Range("I126").Value2 = 1499.75
$ws.Range("K126").Value2 = 4499.25
$ws.Range("M126").Value2 = -2029.25
# WVR row 132
$ws.Range("H132").Value2 = 3980.2104
$ws.Range("I132").Value2 = 3779.516
$ws.Range("K132").Value2 = 11338.548
$ws.Range("M132").Value2 = -8808.548000000001
# WVR row 136
$ws.Range("H136").Value2 = 2386.1714
$ws.Range("I136").Value2 = 1054.2858
$ws.Range("K136").Value2 = 3162.8574
$ws.Range("M136").Value2 = -612.8574000000003

